# Update cryptocurrency price (D) and 1h volume change (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.552.49"
$ws.Range("E2").Value = "  -0.08%  "

$ws.Range("D3").Value = "1.826.19"
$ws.Range("E3").Value = "  -0.12%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.62"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.60%  "

$ws.Range("E6").Value = "  +0.03%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5117"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -5.48%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3957"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.53%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08211"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +6.27%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.117"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.46%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "41.74"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.42%  "

$ws.Range("E12").Value = "  -1.29%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.343"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.07%  "

$ws.Range("E14").Value = "  -0.04%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.558"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.09%  "

$ws.Range("D16").Value = "1.821.05"
$ws.Range("E16").Value = "  -0.40%  "

$ws.Range("E17").Value = "  +3.45%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "93.00"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.19%  "

$ws.Range("E19").Value = "  +0.95%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.85"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.34%  "

$ws.Range("E21").Value = "  +0.01%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.100"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.29%  "

$ws.Range("D23").Value = "28.581.44"
$ws.Range("E23").Value = "  -0.04%  "

$ws.Range("E24").Value = "  +2.05%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.265"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.17%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "21.41"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.83%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "156.68"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.16%  "

$ws.Range("D28").Value = "2.033.31"
$ws.Range("E28").Value = "  -0.31%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.420"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.09%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "127.17"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.21%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.115"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.41%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1089"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.16%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.783"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.62%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.657"
$ws.Range("D34").Style = "Normal"

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.07067"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -6.59%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.2232"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.11%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.289"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.27%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02357"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.19%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.825"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.13%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6338"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.38%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "11.30"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.84%  "

$ws.Range("E42").Value = "  -0.77%  "

$ws.Range("E43").Value = "  -0.26%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.62"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.67%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5954"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.87%  "

$ws.Range("E46").Value = "  +0.64%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "125.29"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.17%  "

$ws.Range("E48").Value = "  -0.24%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.193"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.46%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06951"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.46%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.086"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.47%  "
